$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contains an Excel Table ("Table2") currently spanning A1:H24.
# We need to add one more data row (row 25) to it, growing the table /
# autoFilter / dimension accordingly, then populate the new row with the
# same values/format pattern as the previous row (row 24), but for the
# new task.

$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# Copy the last existing data row's formatting/values down into the new
# row first (this preserves cell types exactly, e.g. the text-stored
# "End date" value), then overwrite the cells that actually change.
$ws.Range("A24:G24").Copy($ws.Range("A25:G25"))

$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "Get DSA Code of XMED from Niha's PC vah vah vivah 4"
$ws.Range("C25").Value = "Completed"
$ws.Range("D25").Value = 100
$ws.Range("G25").Value = "Aishwarrya VP"
# E25 (Start date) and F25 (End date) stay identical to row 24's values,
# already copied above.

$ws.Range("B25").Select()
